$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Map of row number -> new pct_diff value (column G)
$values = @{
    2 = -2.46
    3 = -3.38
    4 = -2.63
    5 = -2.16
    6 = -2.67
    7 = -1.71
    8 = -2.05
    9 = -3.56
    10 = -2.8
    11 = -3.44
    12 = -3.75
    13 = -3.96
    14 = -3.76
    15 = -3.72
    16 = -3.46
    17 = -2.73
    18 = -3.33
    19 = -5.03
    20 = -3.87
    21 = -3.02
    22 = -3.4
    23 = -4.32
    24 = -3.46
    25 = -4.89
    26 = -4.7
    27 = -4.58
    28 = -5.91
    29 = -6.47
    30 = -7.61
    31 = -8.75
    32 = -12.23
    33 = -16.7
    34 = -18.22
    35 = -27.19
    36 = -28.98
    37 = -35.91
    38 = -38.79
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item([int]$row, 7).Value = $values[$row]
}
